$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 344627.44
$ws.Range("J17").Value = 351516.75
$ws.Range("L17").Value = 1054550.25
$ws.Range("N17").Value = -1054886.25
$ws.Range("H62").Value = 1908.909
$ws.Range("I62").Value = 1699.75
$ws.Range("K62").Value = 1699.75
$ws.Range("M62").Value = -1075.75
$ws.Range("H64").Value = 4121.1465
$ws.Range("I64").Value = 4353.148
$ws.Range("K64").Value = 4353.148
$ws.Range("M64").Value = -4105.148
$ws.Range("H65").Value = 1908.909
$ws.Range("I65").Value = 1699.75
$ws.Range("K65").Value = 8498.75
$ws.Range("M65").Value = -5378.75
$ws.Range("H67").Value = 4121.1465
$ws.Range("I67").Value = 4353.148
$ws.Range("K67").Value = 4353.148
$ws.Range("M67").Value = -3495.148
$ws.Range("H129").Value = 937.0925999999999
$ws.Range("J129").Value = 968.60785
$ws.Range("L129").Value = 2905.82355
$ws.Range("N129").Value = -12905.82355
$ws.Range("H137").Value = 2038.1515
$ws.Range("I137").Value = 1455.8334
$ws.Range("J137").Value = 2370.9048
$ws.Range("K137").Value = 4367.5002
$ws.Range("L137").Value = 7112.714399999999
$ws.Range("M137").Value = -1817.5002
$ws.Range("N137").Value = -12212.7144

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 8356.951999999999
$ws.Range("I45").Value = 10947.134
$ws.Range("K45").Value = 10947.134
$ws.Range("M45").Value = -10570.134
$ws.Range("H122").Value = 1605881
$ws.Range("I122").Value = 1712739.8
$ws.Range("K122").Value = 5138219.4
$ws.Range("M122").Value = -5135769.4
$ws.Range("H132").Value = 7795.6
$ws.Range("I132").Value = 10000
$ws.Range("J132").Value = 7244.5
$ws.Range("K132").Value = 30000
$ws.Range("L132").Value = 21733.5
$ws.Range("M132").Value = -27470
$ws.Range("N132").Value = -26793.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H69").Value = 28475
$ws.Range("J69").Value = 28475
$ws.Range("L69").Value = 28475
$ws.Range("N69").Value = -30097
$ws.Range("H72").Value = 28475
$ws.Range("J72").Value = 28475
$ws.Range("L72").Value = 85425
$ws.Range("N72").Value = -93537
$ws.Range("H86").Value = 2425.0625
$ws.Range("I86").Value = 2300.0715
$ws.Range("J86").Value = 3300
$ws.Range("K86").Value = 2300.0715
$ws.Range("L86").Value = 3300
$ws.Range("M86").Value = -1177.0715
$ws.Range("N86").Value = -5546
$ws.Range("H89").Value = 2425.0625
$ws.Range("I89").Value = 2300.0715
$ws.Range("J89").Value = 3300
$ws.Range("K89").Value = 11500.3575
$ws.Range("L89").Value = 16500
$ws.Range("M89").Value = -5884.3575
$ws.Range("N89").Value = -27732
$ws.Range("H138").Value = 59745.715
$ws.Range("J138").Value = 59745.715
$ws.Range("L138").Value = 59745.715
$ws.Range("N138").Value = -70025.715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3723.7087
$ws.Range("I31").Value = 1888.7142
$ws.Range("J31").Value = 4388.1035
$ws.Range("K31").Value = 1888.7142
$ws.Range("L31").Value = 4388.1035
$ws.Range("M31").Value = -1593.7142
$ws.Range("N31").Value = -4978.1035
$ws.Range("H34").Value = 3723.7087
$ws.Range("I34").Value = 1888.7142
$ws.Range("J34").Value = 4388.1035
$ws.Range("K34").Value = 1888.7142
$ws.Range("L34").Value = 4388.1035
$ws.Range("M34").Value = -1686.7142
$ws.Range("N34").Value = -4792.1035
$ws.Range("H63").Value = 41299.5
$ws.Range("J63").Value = 41299.5
$ws.Range("L63").Value = 41299.5
$ws.Range("N63").Value = -42671.5
$ws.Range("H66").Value = 41299.5
$ws.Range("J66").Value = 41299.5
$ws.Range("L66").Value = 123898.5
$ws.Range("N66").Value = -130762.5
$ws.Range("H132").Value = 4259.8945
$ws.Range("I132").Value = 3010
$ws.Range("K132").Value = 9030
$ws.Range("M132").Value = -6500

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 4559.2
$ws.Range("I63").Value = 1664
$ws.Range("J63").Value = 5800
$ws.Range("K63").Value = 4992
$ws.Range("L63").Value = 17400
$ws.Range("M63").Value = -4243
$ws.Range("N63").Value = -18898
$ws.Range("H66").Value = 4559.2
$ws.Range("I66").Value = 1664
$ws.Range("J66").Value = 5800
$ws.Range("K66").Value = 14976
$ws.Range("L66").Value = 52200
$ws.Range("M66").Value = -11232
$ws.Range("N66").Value = -59688
$ws.Range("H68").Value = 2682.923
$ws.Range("I68").Value = 3927.0303
$ws.Range("J68").Value = 1770.5778
$ws.Range("K68").Value = 11781.0909
$ws.Range("L68").Value = 5311.7334
$ws.Range("M68").Value = -10970.0909
$ws.Range("N68").Value = -6933.7334
$ws.Range("H71").Value = 2682.923
$ws.Range("I71").Value = 3927.0303
$ws.Range("J71").Value = 1770.5778
$ws.Range("K71").Value = 35343.2727
$ws.Range("L71").Value = 15935.2002
$ws.Range("M71").Value = -31287.2727
$ws.Range("N71").Value = -24047.2002
$ws.Range("H112").Value = 2743.9614
$ws.Range("J112").Value = 3060.8696
$ws.Range("L112").Value = 9182.6088
$ws.Range("N112").Value = -11398.6088

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5866.5747
$ws.Range("I70").Value = 5937.9414
$ws.Range("J70").Value = 5679.923
$ws.Range("K70").Value = 5937.9414
$ws.Range("L70").Value = 5679.923
$ws.Range("M70").Value = -5667.9414
$ws.Range("N70").Value = -6219.923
$ws.Range("H73").Value = 5866.5747
$ws.Range("I73").Value = 5937.9414
$ws.Range("J73").Value = 5679.923
$ws.Range("K73").Value = 5937.9414
$ws.Range("L73").Value = 5679.923
$ws.Range("M73").Value = -5001.9414
$ws.Range("N73").Value = -7551.923
$ws.Range("H122").Value = 7656287
$ws.Range("I122").Value = 4987867.5
$ws.Range("J122").Value = 25001012
$ws.Range("K122").Value = 14963602.5
$ws.Range("L122").Value = 75003036
$ws.Range("M122").Value = -14961152.5
$ws.Range("N122").Value = -75007936
$ws.Range("H132").Value = 4844.913
$ws.Range("I132").Value = 16750
$ws.Range("J132").Value = 3711.0952
$ws.Range("K132").Value = 50250
$ws.Range("L132").Value = 11133.2856
$ws.Range("M132").Value = -47720
$ws.Range("N132").Value = -16193.2856

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 111111656
$ws.Range("I107").Value = 166667150
$ws.Range("J107").Value = 667
$ws.Range("K107").Value = 500001450
$ws.Range("L107").Value = 2001
$ws.Range("M107").Value = -499999530
$ws.Range("N107").Value = -5841
